$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "TP08" column (J) that duplicates column I's header/formatting.
# Column I already holds the "TP08" results; J is a fresh copy of the same
# grading column (same header + same cell borders/shading), but only the
# first 16 students (rows 2-17) already have a "1" filled in - the rest
# (rows 18-48) are left blank, ready to be filled in later.

# 1) Copy all formatting (styles/borders/fills) from I1:I48 down to J1:J48
$ws.Range("I1:I48").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats

# 2) Copy the values (header text + the "1" marks) from I1:I17 into J1:J17
$ws.Range("I1:I17").Copy()
$ws.Range("J1").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0

# 3) Update the selected/active cell to reflect where the author left off
[void]$ws.Range("J17").Select()
